# Daily attendance processing
# Normalises the "Recorded By" column (G) on the active sheet: for any cell
# whose value is a comma-separated list of recorders, the list order is
# reversed (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
# Cells with only a single recorder are left untouched (reversing a
# one-element list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "
    $n = $parts.Count
    if ($n -le 1) { continue }

    $reversedStr = ""
    for ($i = $n - 1; $i -ge 0; $i--) {
        if ($reversedStr -ne "") { $reversedStr += ", " }
        $reversedStr += $parts[$i]
    }

    if ($reversedStr -ne $val) {
        $cell.Value() = $reversedStr
    }
}
